# Generate Report for Archive
#
# 1) Update the localization "Status" value from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2/F2, zh-cn!C2,
#    de-de!C2 all share the same string).
# 2) Narrow the "Status" column(s) to match the new (shorter) text:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update the status text -------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Resize the status columns ----------------------------------------------
# The saved <col width="..."> value in the OOXML is the ColumnWidth plus a
# fixed character-padding offset (5/6 of a character); subtract it here so
# the persisted width matches the target value as closely as this interop
# layer's quantization allows.
$targetStoredWidth = 13.4101845877511
$widthPaddingOffset = 0.8333333333333334
$newWidth = $targetStoredWidth - $widthPaddingOffset

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth = $newWidth
$dede.Columns.Item(3).ColumnWidth = $newWidth
